$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2025-04-28 08:01:43"
$ws.Range("B7").Value = 204
